$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header C1: "Nmbre y Apellido" -> "Nombre y Apellido"
$ws.Range("C1").Value = "Nombre y Apellido"

# Row 2: replace the previous winner (112 / Marciana Garay.) with the new one (14 / Sergio Riquelme.)
$ws.Range("A2").Value = 14
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "14"
$ws.Range("C2").Value = "Sergio Riquelme."
$ws.Range("D2").Value = "si"
$ws.Range("E2").Value = "Ganador de Gs. 1.000.000"

# Row 3 (new): 18 / Jorge Morinigo.
$ws.Range("A3").Value = 18
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "18"
$ws.Range("C3").Value = "Jorge Morinigo."
$ws.Range("D3").Value = "si"
$ws.Range("E3").Value = "Ganador de Gs. 1.000.000"

# Row 4 (new): 74 / Stella Martinez.
$ws.Range("A4").Value = 74
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "74"
$ws.Range("C4").Value = "Stella Martinez."
$ws.Range("D4").Value = "si"
$ws.Range("E4").Value = "Ganador de Gs. 1.000.000"
